$d = $word.ActiveDocument
$full = $d.Content

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body><w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">To</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Christ</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">To</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">the</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Land</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr><w:r><w:t xml:space="preserve">Dorothy</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Day</w:t></w:r></w:p><w:p><w:r><w:rPr><w:i /></w:rPr><w:t xml:space="preserve">The Catholic Worker</w:t></w:r><w:r><w:t xml:space="preserve">, January 1936, 1-2.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:i /></w:rPr><w:t xml:space="preserve">Summary: Presents P. Maurin three-point program: Round Table Discussions, Houses of Hospitality, and Farming Communes to further the personalist and communitarian revolution. Promotes worker ownership in order to go back to the land to establish farming communes. (DDLW #143).</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">For those who have put to us the question &quot;What have you to offer in the way of a constructive program for a new social order?&quot; we have replied over and over, &quot;Peter Maurin&#39;s three-point program</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i /></w:rPr><w:t xml:space="preserve">of Round-table Discussions, Houses of Hospitality, Farming Communes</w:t></w:r><w:r><w:t xml:space="preserve">.&quot; This program is so simple as to be unsatisfactory to most, who look for something to be complicated before it can be successful. Remembering the words of St. Francis that we cannot know what we have not practiced, we have tried not only to publish a paper but to put our program into practice. From the very beginning we have sought clarification of thought through</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i /></w:rPr><w:t xml:space="preserve">The Catholic Worker</w:t></w:r><w:r><w:t xml:space="preserve">, through round-table discussions, forums, through circulating literature. We have had a workers&#39; school where the finest scholars of the Church have come to teach. We have had a House of Hospitality now for two years, where we gave shelter to the homeless, fed the hungry, clothed the naked, and cared for the sick. We have tried, all of us, to be workers and scholars, and to combine work and prayer according to the Benedictine ideal. We have tried to imitate St. Francis in his holy poverty. Our aim has been to combat the atheism of the day by our devotion to the liturgical movement; to combat the bourgeois spirit by the Franciscan spirit; to oppose to class-war technique the performance of the works of mercy.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">We have not altogether neglected the farming commune idea, inasmuch as we had a halfway house in Staten Island where children were given vacations, weekend conferences were held and the sick cared for, and a garden cultivated.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">March 1 will see the start of a serious attempt to put into practice the third point of our program. We are going to move out on a farm, within a few hours of New York, and start there a true farming commune.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">We are making this move because we do not feel that we can talk in the paper about something we are not practicing. We believe that our words will have more weight, our writings will carry more conviction, I we ourselves are engaged in making a better life on the land.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">This does not mean that we are going to abandon the city, which we realize is above all the home of the dispossessed, of the forgotten. We shall keep a group in New York City and the work of the apostolate of labor will go on. We shall also be sending out apostles of labor from the farm, to scenes of industrial conflict, to factories and to lodging houses to live and work with the poor. The columns of the paper will be filled as usual with industrial news, discussion of unionism, the cooperative movement, maternity guilds, relief, public and private. But there will be more space devoted to rural life problems, and you will hear from month to month how the work of the farming commune is progressing, the difficulties, the mistakes, and the progress of the work.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Help us in this venture, which is your venture, too. And pray with us that we get out of the city by March 1.</w:t></w:r></w:p></w:body>' + `
'</w:document>' + `
'</pkg:xmlData></pkg:part></pkg:package>'

$full.InsertXML($xml)
$end = $d.Content.End
Write-Host "Done. Content end:" $end
for ($off = -5; $off -le 0; $off++) {
    $pos = $end + $off
    if ($pos -ge 0) {
        try {
            $t = $d.Range($pos, $end).Text
            Write-Host "range($pos,$end) text=[$t]"
        } catch {
            Write-Host "pos $pos ERR $_"
        }
    }
}
